# Update "Pais" (countries) data sheet with refreshed COVID stats and
# re-ordered country rankings (some countries swapped ranking position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 01:42"

# --- Country name swaps caused by re-ranking (same row positions, new order) ---
# Rows 50/51: Rumania <-> Nigeria
$ws.Range("A50").Value = "Nigeria"
$ws.Range("A51").Value = "Rumania"

# Rows 83/84: Consejo Danes para los Refugiados <-> Etiopia
$ws.Range("A83").Value = "Etiopia"
$ws.Range("A84").Value = "Consejo Danes para los Refugiados"

# Rows 142/143: Liberia <-> Republica de Chipre
$ws.Range("A142").Value = "Republica de Chipre"
$ws.Range("A143").Value = "Liberia"

# Rows 209/210: Groenlandia <-> Islas Malvinas
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Updated numeric statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3614461
$ws.Range("C4").Value = 69384
$ws.Range("D4").Value = 1637726
$ws.Range("E4").Value = 1836685
$ws.Range("G4").Value = 907
$ws.Range("H4").Value = 140050

# Row 5: Brasil
$ws.Range("B5").Value = 1970909
$ws.Range("C5").Value = 39705
$ws.Range("D5").Value = 1255564
$ws.Range("E5").Value = 639822
$ws.Range("G5").Value = 1261
$ws.Range("H5").Value = 75523

# Row 31: Ecuador
$ws.Range("B31").Value = 70329
$ws.Range("C31").Value = 759
$ws.Range("D31").Value = 30641
$ws.Range("E31").Value = 34530
$ws.Range("G31").Value = 28
$ws.Range("H31").Value = 5158

# Row 50: Nigeria (after swap)
$ws.Range("B50").Value = 34259
$ws.Range("C50").Value = 643
$ws.Range("D50").Value = 13999
$ws.Range("E50").Value = 19500
$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 760

# Row 51: Rumania (after swap)
$ws.Range("B51").Value = 34226
$ws.Range("C51").Value = 641
$ws.Range("D51").Value = 22049
$ws.Range("E51").Value = 10225
$ws.Range("G51").Value = 21
$ws.Range("H51").Value = 1952

# Row 54: Guatemala
$ws.Range("B54").Value = 32074
$ws.Range("C54").Value = 1202
$ws.Range("D54").Value = 4624
$ws.Range("E54").Value = 26100
$ws.Range("G54").Value = 48
$ws.Range("H54").Value = 1350

# Row 58: Ghana
$ws.Range("B58").Value = 25430
$ws.Range("C58").Value = 442
$ws.Range("D58").Value = 21511
$ws.Range("E58").Value = 3780

# Row 59: Japon
$ws.Range("B59").Value = 22508
$ws.Range("C59").Value = 288
$ws.Range("D59").Value = 18545
$ws.Range("E59").Value = 2979
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 984

# Row 69: Chequia
$ws.Range("B69").Value = 13475
$ws.Range("C69").Value = 134
$ws.Range("D69").Value = 8507
$ws.Range("E69").Value = 4613

# Row 78: Noruega
$ws.Range("B78").Value = 9011
$ws.Range("C78").Value = 10
$ws.Range("E78").Value = 620

# Row 83: Etiopia (after swap)
$ws.Range("B83").Value = 8181
$ws.Range("C83").Value = 212
$ws.Range("D83").Value = 2430
$ws.Range("E83").Value = 5605
$ws.Range("G83").Value = 7
$ws.Range("H83").Value = 146

# Row 84: Consejo Danes para los Refugiados (after swap)
$ws.Range("B84").Value = 8163
$ws.Range("C84").Value = 28
$ws.Range("D84").Value = 3983
$ws.Range("E84").Value = 3988
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 192

# Row 92: Guinea
$ws.Range("B92").Value = 6276
$ws.Range("C92").Value = 76
$ws.Range("D92").Value = 4981
$ws.Range("E92").Value = 1257

# Row 93: Gabon
$ws.Range("B93").Value = 6121
$ws.Range("C93").Value = 95
$ws.Range("D93").Value = 3664
$ws.Range("E93").Value = 2411

# Row 98: Republica de Africa Central
$ws.Range("B98").Value = 4362
$ws.Range("C98").Value = 6
$ws.Range("D98").Value = 1261
$ws.Range("E98").Value = 3048

# Row 106: Somalia
$ws.Range("B106").Value = 3083
$ws.Range("C106").Value = 7
$ws.Range("D106").Value = 1425
$ws.Range("E106").Value = 1565

# Row 109: Mayotte
$ws.Range("B109").Value = 2743
$ws.Range("C109").Value = 6
$ws.Range("D109").Value = 2581
$ws.Range("E109").Value = 125

# Row 142: Republica de Chipre (after swap)
$ws.Range("B142").Value = 1025
$ws.Range("C142").Value = 2
$ws.Range("D142").Value = 839
$ws.Range("E142").Value = 167
$ws.Range("H142").Value = 19

# Row 143: Liberia (after swap)
$ws.Range("B143").Value = 1024
$ws.Range("D143").Value = 439
$ws.Range("E143").Value = 534
$ws.Range("H143").Value = 51

# Row 149: Surinam
$ws.Range("B149").Value = 837
$ws.Range("C149").Value = 36
$ws.Range("E149").Value = 246

# Row 181: Bahamas
$ws.Range("B181").Value = 119
$ws.Range("C181").Value = 6
$ws.Range("E181").Value = 17

# Row 184: Barbados
$ws.Range("B184").Value = 104
$ws.Range("C184").Value = 1
$ws.Range("E184").Value = 7
